{"js": "// \"formating daftar gambar, isi, lampiran, tabel, kata pengantar\"\n//\n// The \"DAFTAR LAMPIRAN\" title paragraph had an explicit run-level / paragraph\n// -mark font size override (sz=24 / szCs=24, i.e. 12pt) that was fighting the\n// Heading1 style's own size (sz=28 / szCs=32, i.e. 14pt/16pt). The fix removes\n// the explicit <w:sz> override (so the heading falls back to the style's\n// 14pt) and bumps the complex-script size to szCs=28 (14pt) to match, and\n// also moves the title run to sit *inside* the _GoBack bookmark span instead\n// of after it.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"DAFTAR LAMPIRAN\" title paragraph (first paragraph of the doc).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"DAFTAR LAMPIRAN\") {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"DAFTAR LAMPIRAN heading paragraph not found\");\n}\n\n// Rebuild the paragraph's OOXML: drop the explicit <w:sz> (let Heading1's own\n// size win), set szCs to 28 (14pt) on both the paragraph mark and the run,\n// and place the title run between the bookmarkStart/bookmarkEnd pair.\nconst newXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  \"<w:pPr>\" +\n  '<w:pStyle w:val=\"Heading1\"/>' +\n  '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"0\"/></w:numPr>' +\n  '<w:spacing w:before=\"0\" w:after=\"0\" w:line=\"480\" w:lineRule=\"auto\"/>' +\n  '<w:jc w:val=\"center\"/>' +\n  '<w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr>' +\n  \"</w:pPr>\" +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  \"<w:r>\" +\n  '<w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr>' +\n  \"<w:t>DAFTAR LAMPIRAN</w:t>\" +\n  \"</w:r>\" +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst range = target.getRange();\nrange.insertOoxml(newXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"formating daftar gambar, isi, lampiran, tabel, kata pengantar\"\n#\n# The \"DAFTAR LAMPIRAN\" title paragraph had an explicit run-level / paragraph-\n# mark font size override (sz=24 / szCs=24, i.e. 12pt) that was fighting the\n# Heading1 style's own size (sz=28 / szCs=32, i.e. 14pt/16pt). The fix removes\n# the explicit <w:sz> override (so the heading falls back to the style's\n# 14pt) and bumps the complex-script size to szCs=28 (14pt) to match, and\n# also moves the title run to sit *inside* the _GoBack bookmark span instead\n# of after it.\n\n$d = $word.ActiveDocument\n\n# Locate the \"DAFTAR LAMPIRAN\" title paragraph (first paragraph of the doc).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text -replace \"[\\r\\a]+$\", \"\"\n    if ($t -eq \"DAFTAR LAMPIRAN\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"DAFTAR LAMPIRAN heading paragraph not found\"\n}\n\n# Rebuild the paragraph's OOXML: drop the explicit <w:sz> (let Heading1's own\n# size win), set szCs to 28 (14pt) on both the paragraph mark and the run,\n# and place the title run between the bookmarkStart/bookmarkEnd pair.\n$newXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:pPr>\n<w:pStyle w:val=\"Heading1\"/>\n<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"0\"/></w:numPr>\n<w:spacing w:before=\"0\" w:after=\"0\" w:line=\"480\" w:lineRule=\"auto\"/>\n<w:jc w:val=\"center\"/>\n<w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr>\n</w:pPr>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:r>\n<w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr>\n<w:t>DAFTAR LAMPIRAN</w:t>\n</w:r>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$target.Range.InsertXML($newXml)\n"}
